$d = $word.ActiveDocument

# --- 1) "Made a github repository for the app" paragraph: collapse the
#        spell-check-split runs (and proofErr markers) into a single run.
$p8 = $d.Paragraphs.Item(8)
$r8 = $p8.Range
$r8.Find.Execute("Made a github repository for the app", $true, $false, $false, $false, $false, $true, 1, $false, "Made a github repository for the app", 2)

# --- 2) "First git commit and push to the app repo" paragraph: collapse the
#        grammar-check-split runs (and proofErr markers) into a single run.
$p10 = $d.Paragraphs.Item(10)
$r10 = $p10.Range
$r10.Find.Execute("First git commit and push to the app repo", $true, $false, $false, $false, $false, $true, 1, $false, "First git commit and push to the app repo", 2)

# --- 3) Insert a brand new paragraph right after it, made up of three runs.
$newParAnchor = $d.Paragraphs.Item(10).Range
$newParAnchor.InsertParagraphAfter()

$newPar = $d.Paragraphs.Item(11).Range
$newPar.InsertAfter("Had a meeting with all members, explained working of app to vishal informed him his work, ")

$start2 = $d.Paragraphs.Item(11).Range.End - 1
$run2Text = "taught him the git push commands, tested 2 commit and push"
$d.Range($start2, $start2).InsertAfter($run2Text)

$start3 = $d.Paragraphs.Item(11).Range.End - 1
$run3Text = ", decided to input the plant data in an excel sheet, plant data has name of plant, one youtube to explain growing of the plant and development of seed, one document for reading about the plant."
$d.Range($start3, $start3).InsertAfter($run3Text)

# Force the 2nd and 3rd runs to be kept as distinct <w:r> elements (instead of
# being silently coalesced back into the run before them) by touching a
# direct-formatting property and reverting it right away. Order matters:
# do the later run first, then the earlier one.
$range3 = $d.Range($start3, $start3 + $run3Text.Length)
$range3.Font.Bold = 1
$range3.Font.Bold = 0

$range2 = $d.Range($start2, $start2 + $run2Text.Length)
$range2.Font.Bold = 1
$range2.Font.Bold = 0

Write-Output "done"
